# Edit of the PPT
#
# 1) Slide 3 ("Sommaire"): the two runs "Présentation " + "projet" in the
#    content placeholder get merged into a single run "Présentation projet".
# 2) Slides 2-7 (all but the title slide) get the "Slide Number" placeholder
#    turned on, which inserts a sldNum placeholder shape copied from the
#    layout (the shape's name follows PowerPoint's usual
#    "Espace réservé du numéro de diapositive N" pattern).

$p = $ppt.ActivePresentation

# --- 1) Merge the "Présentation " / "projet" runs on slide 3 -----------
$s3 = $p.Slides.Item(3)
$contentShape = $s3.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
# Force the engine to actually rewrite the run (a no-op same-text
# assignment is short-circuited), then set the final merged text.
$para.Text = "__tmp__"
$para2 = $tr.Paragraphs(3, 1)
$para2.Text = "Présentation projet"

# --- 2) Turn on slide numbers for every slide except the title slide ---
$slideNumberNames = @{
    2 = "Espace réservé du numéro de diapositive 3"
    3 = "Espace réservé du numéro de diapositive 3"
    4 = "Espace réservé du numéro de diapositive 2"
    5 = "Espace réservé du numéro de diapositive 3"
    6 = "Espace réservé du numéro de diapositive 3"
    7 = "Espace réservé du numéro de diapositive 3"
}

foreach ($idx in 2, 3, 4, 5, 6, 7) {
    $s = $p.Slides.Item($idx)
    $s.HeadersFooters.SlideNumber.Visible = $true
    $newShape = $s.Shapes.Item($s.Shapes.Count)
    $newShape.Name = $slideNumberNames[$idx]
}
